$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a blank row 2 separating the header row (row 1) from the
# data row (originally row 3). Delete the blank row so the data row moves
# up to become row 2.
$ws.Rows(2).Delete()

# The surviving rows were carrying custom (bold/centered/right-aligned,
# larger) fonts. Strip all that explicit formatting back to the workbook
# default.
$ws.Range("A1:D2").ClearFormats()

# The workbook's default/"Normal" font also changes, from the old
# Arial-based header font to 10pt Calibri.
$normal = $wb.Styles.Item("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.Size = 10

# Keep the active selection in sync with the row shift (was D3, now D2).
[void]$ws.Range("D2").Select()
